$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

$ws.Range("E3").Value = "v2024-04-30"
$ws.Range("E4").Value = "v3.65.0"

$ws.Range("E4").Select()
